{"js": "// 1) Fix the spell-checked typo \"quorium\" -> \"quorum\" and merge the\n//    \"D\u00e9tection du \" / \"quorium\" runs back into a single\n//    \"D\u00e9tection du quorum\" run (removing the split caused by the\n//    now-resolved spelling flag).\nconst quoriumResults = context.document.body.search(\"quorium\", { matchCase: true });\nquoriumResults.load(\"items\");\nawait context.sync();\nquoriumResults.items[0].delete();\nawait context.sync();\n\nconst detectionResults = context.document.body.search(\"D\u00e9tection du \", { matchCase: true });\ndetectionResults.load(\"items\");\nawait context.sync();\ndetectionResults.items[0].insertText(\"quorum\", Word.InsertLocation.end);\nawait context.sync();\n\n// 2) Add the missing comma: \"st\u00e9ro\u00efdes vitamine D\" -> \"st\u00e9ro\u00efdes, vitamine D\"\nconst mediatorResults = context.document.body.search(\"st\u00e9ro\u00efdes vitamine D\", { matchCase: true });\nmediatorResults.load(\"items\");\nawait context.sync();\nmediatorResults.items[0].insertText(\"st\u00e9ro\u00efdes, vitamine D\", Word.InsertLocation.replace);\nawait context.sync();\n\n// 3) Turn the empty paragraph right after the last table into a new\n//    remark: an underlined \"Rmq :\" label followed by the remark text.\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\nconst lastParagraph = paragraphs.items[paragraphs.items.length - 1];\n\n// Insert the (plain) remark text first so it lands in its own run,\n// then insert the \"Rmq :\" label before it and underline just that run.\nconst remarkRun = lastParagraph.insertText(\n  \" les mol\u00e9cules liposolubles p\u00e9n\u00e8trent dans toutes les cellules mais une r\u00e9ponse cellulaire ne se produira que si la cellule est dot\u00e9e des r\u00e9cepteurs ad\u00e9quates. \",\n  Word.InsertLocation.start\n);\nawait context.sync();\n\nconst labelRun = lastParagraph.insertText(\"Rmq :\", Word.InsertLocation.start);\nawait context.sync();\nlabelRun.font.underline = Word.UnderlineType.single;\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# 1) Fix the spell-checked typo \"quorium\" -> \"quorum\" and merge the\n#    \"D\u00e9tection du \" / \"quorium\" runs back into a single\n#    \"D\u00e9tection du quorum\" run (removing the split caused by the\n#    now-resolved spelling flag).\n$quoriumRange = $d.Content\n$quoriumFind = $quoriumRange.Find\n$quoriumFind.Text = \"quorium\"\n$quoriumFind.Execute() | Out-Null\n$quoriumRange.Delete()\n\n$detectionRange = $d.Content\n$detectionFind = $detectionRange.Find\n$detectionFind.Text = \"D\u00e9tection du \"\n$detectionFind.Execute() | Out-Null\n$detectionRange.Collapse(0)  # wdCollapseEnd\n$detectionRange.Text = \"quorum\"\n\n# 2) Add the missing comma: \"st\u00e9ro\u00efdes vitamine D\" -> \"st\u00e9ro\u00efdes, vitamine D\"\n$mediatorRange = $d.Content\n$mediatorFind = $mediatorRange.Find\n$mediatorFind.Text = \"st\u00e9ro\u00efdes vitamine D\"\n$mediatorFind.Execute() | Out-Null\n$mediatorRange.Text = \"st\u00e9ro\u00efdes, vitamine D\"\n\n# 3) Turn the empty paragraph right after the last table into a new\n#    remark: an underlined \"Rmq :\" label followed by the remark text.\n$lastParagraph = $d.Paragraphs.Last\n\n$remarkRange = $lastParagraph.Range\n$remarkRange.Collapse(1)  # wdCollapseStart\n$remarkRange.InsertAfter(\" les mol\u00e9cules liposolubles p\u00e9n\u00e8trent dans toutes les cellules mais une r\u00e9ponse cellulaire ne se produira que si la cellule est dot\u00e9e des r\u00e9cepteurs ad\u00e9quates. \")\n\n$labelRange = $lastParagraph.Range\n$labelRange.Collapse(1)  # wdCollapseStart\n$labelRange.InsertAfter(\"Rmq :\")\n$labelRange.Font.Underline = 1  # wdUnderlineSingle\n"}
